$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.909.14'
$ws.Range("E2").Value = '  -1.40%  '

$ws.Range("D3").Value = '1.636.81'
$ws.Range("E3").Value = '  -0.65%  '

$ws.Range("D5").Value = '215.51'
$ws.Range("E5").Value = '  -0.64%  '

$ws.Range("E6").Value = '  +0.25%  '

$ws.Range("E9").Value = '  -0.27%  '

$ws.Range("D10").Value = '19.55'
$ws.Range("E10").Value = '  -1.95%  '

$ws.Range("D11").Value = '0.0793'
$ws.Range("E11").Value = '  -0.08%  '

$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = '4.28'

$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '1.863.29'
$ws.Range("E13").Value = '  -0.62%  '

$ws.Range("D14").Value = '1.654.12'
$ws.Range("E14").Value = '  +0.02%  '

$ws.Range("E15").Value = '  -0.81%  '

$ws.Range("D16").Value = '0.0₃0764'
$ws.Range("E16").Value = '  -0.36%  '

$ws.Range("D17").Value = '62.82'
$ws.Range("E17").Value = '  -0.87%  '

$ws.Range("D18").Value = '25.966.85'
$ws.Range("E18").Value = '  -1.15%  '

$ws.Range("E19").Value = '  +0.18%  '

$ws.Range("D20").Value = '192.93'
$ws.Range("E20").Value = '  -1.48%  '

$ws.Range("D21").Value = '4.36'
$ws.Range("E21").Value = '  -2.12%  '

$ws.Range("D22").Value = '9.91'
$ws.Range("E22").Value = '  -1.74%  '

$ws.Range("E23").Value = '  -1.05%  '

$ws.Range("D24").Value = '0.131'
$ws.Range("E24").Value = '  +4.37%  '

$ws.Range("E25").Value = '  -0.04%  '

$ws.Range("D26").Value = '143.22'
$ws.Range("E26").Value = '  -0.16%  '

$ws.Range("E27").Value = '  +0.14%  '

$ws.Range("E28").Value = '  -1.01%  '

$ws.Range("D29").Value = '15.55'

$ws.Range("E30").Value = '  -0.68%  '

$ws.Range("E31").Value = '  -0.82%  '

$ws.Range("E32").Value = '  -2.23%  '

$ws.Range("D33").Value = '3.24'
$ws.Range("E33").Value = '  -0.43%  '

$ws.Range("D34").Value = '1.54'
$ws.Range("E34").Value = '  -4.28%  '

$ws.Range("E35").Value = '  +1.61%  '

$ws.Range("D36").Value = '0.900'
$ws.Range("E36").Value = '  -1.52%  '

$ws.Range("D37").Value = '1.133.54'
$ws.Range("E37").Value = '  -0.51%  '

$ws.Range("D38").Value = '0.542'
$ws.Range("E38").Value = '  -2.23%  '

$ws.Range("E39").Value = '  -1.48%  '

$ws.Range("E40").Value = '  -0.52%  '

$ws.Range("E41").Value = '  -0.81%  '

$ws.Range("D42").Value = '99.25'
$ws.Range("E42").Value = '  -1.29%  '

$ws.Range("E43").Value = '  -0.42%  '

$ws.Range("D44").Value = '1.773.07'
$ws.Range("E44").Value = '  -0.57%  '

$ws.Range("D45").Value = '0.0₆0115'
$ws.Range("E45").Value = '  +2.19%  '

$ws.Range("D46").Value = '56.58'
$ws.Range("E46").Value = '  -0.98%  '

$ws.Range("D47").Value = '0.0528'
$ws.Range("E47").Value = '  +2.14%  '

$ws.Range("E48").Value = '  -1.01%  '

$ws.Range("D49").Value = '7.70'
$ws.Range("E49").Value = '  -0.21%  '

$ws.Range("E50").Value = '  -0.88%  '

$ws.Range("D51").Value = '0.0959'
$ws.Range("E51").Value = '  -1.33%  '
